# Weekly price-sheet update: insert one new daily quote row for
# "Poroto granado" (Macroferia Regional de Talca) ahead of the existing
# row 108, shifting all later rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 108 (pushes old rows 108..135 down to 109..136,
# and Excel copies the row-above formatting, e.g. the date style in column D).
$ws.Rows.Item(108).Insert()

# Populate the new row with the latest weekly quote.
$ws.Range("A108").Value = 5
$ws.Range("B108").Value = "Macroferia Regional de Talca"
$ws.Range("C108").Value = "Maule"
$ws.Range("D108").Value = 44642
$ws.Range("E108").Value = 7
$ws.Range("F108").Value = 100112030
$ws.Range("G108").Value = "Poroto granado"
$ws.Range("H108").Value = "Sin especificar"
$ws.Range("I108").Value = "Primera"
$ws.Range("J108").Value = 300
$ws.Range("K108").Value = 20000
$ws.Range("L108").Value = 20000
$ws.Range("M108").Value = 20000
$ws.Range("N108").Value = "$/saco 25 kilos"
$ws.Range("O108").Value = "Región del Maule"
$ws.Range("P108").Value = 800
$ws.Range("Q108").Value = 25
$ws.Range("R108").Value = "Hortaliza"
